# Add merge cells (B13:E13 and B22:E22) to each of the existing "Week 1".."Week 4"
# sheets -- these ranges already hold the "Wednesday" / "Friday" day headers but
# were not merged across B:E like the "Monday" header (B4:E4) already is.
$wb = $excel.ActiveWorkbook

for ($i = 1; $i -le 4; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Range("B13:E13").Merge()
    $ws.Range("B22:E22").Merge()
}

# Create the new "Week 5" sheet by duplicating "Week 4" (so it starts out with the
# same look & feel / styles / merged title-and-day-header cells) and then editing
# its contents in place to match the Week 5 program.
$ws4 = $wb.Worksheets.Item(4)
$ws4.Copy($null, $ws4)
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "Week 5"

# Week 5 has an extra "top single" row under each day's Bench/OHP work compared to
# Week 4, so insert one blank row after each of the three "main lift" rows (6, 15,
# 24) -- inserting bottom-up so earlier row numbers stay valid.
$ws5.Rows.Item(25).Insert()
$ws5.Rows.Item(16).Insert()
$ws5.Rows.Item(7).Insert()

# Title
$ws5.Range("B1").Value = "Week 5: Strength 1 of 3"

# Monday block
$ws5.Range("C6").Value = 1
$ws5.Range("D6").Value = 1
$ws5.Range("E6").Value = 365

$ws5.Range("C7").Value = 3
$ws5.Range("D7").Value = "2RIR"
$ws5.Range("E7").Value = 320

$ws5.Range("D9").Value = "6-10"

# Wednesday block
$ws5.Range("C16").Value = 1
$ws5.Range("D16").Value = 1
$ws5.Range("E16").Value = 375

$ws5.Range("C17").Value = 3
$ws5.Range("D17").Value = "2RIR"
$ws5.Range("E17").Value = 330

$ws5.Range("D19").Value = "6-10"

# Friday block
$ws5.Range("C26").Value = 1
$ws5.Range("D26").Value = 1
$ws5.Range("E26").Value = 270

$ws5.Range("C27").Value = 3
$ws5.Range("D27").Value = "2RIR"
$ws5.Range("E27").Value = 235

$ws5.Range("D29").Value = 8

# Merge the day header rows (Monday/Wednesday/Friday) the same way the existing
# weeks now are; B1/B4 merges were already carried over from the Week 4 copy.
$ws5.Range("B14:E14").Merge()
$ws5.Range("B24:E24").Merge()

# Make "Week 5" the active sheet/tab, matching the saved workbook view state.
$ws5.Activate()
